# Refresh scorecard data (scorecard_2): updated per-hole stats for Jerome.
# Numeric-looking values are written with a leading apostrophe so Excel keeps
# them as text (matching the sheet's existing shared-string / text-as-data style)
# instead of silently coercing them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Distance
$ws.Range('B2').Value = '''483'
$ws.Range('C2').Value = '''320'
$ws.Range('D2').Value = '''163'
$ws.Range('E2').Value = '''330'
$ws.Range('F2').Value = '''476'
$ws.Range('G2').Value = '''304'
$ws.Range('H2').Value = '''312'
$ws.Range('I2').Value = '''151'
$ws.Range('J2').Value = '''341'
$ws.Range('K2').Value = '''2879'
$ws.Range('M2').Value = '''133'
$ws.Range('N2').Value = '''310'
$ws.Range('O2').Value = '''437'
$ws.Range('P2').Value = '''278'
$ws.Range('Q2').Value = '''317'
$ws.Range('R2').Value = '''163'
$ws.Range('S2').Value = '''278'
$ws.Range('T2').Value = '''379'
$ws.Range('U2').Value = '''2764'
$ws.Range('V2').Value = '''5644'

# Row 3 - Handicap
$ws.Range('B3').Value = '''2'
$ws.Range('C3').Value = '''16'
$ws.Range('D3').Value = '''18'
$ws.Range('E3').Value = '''10'
$ws.Range('F3').Value = '''6'
$ws.Range('G3').Value = '''14'
$ws.Range('H3').Value = '''12'
$ws.Range('I3').Value = '''4'
$ws.Range('J3').Value = '''8'
$ws.Range('L3').Value = '''5'
$ws.Range('M3').Value = '''17'
$ws.Range('N3').Value = '''15'
$ws.Range('O3').Value = '''3'
$ws.Range('P3').Value = '''13'
$ws.Range('Q3').Value = '''9'
$ws.Range('R3').Value = '''7'
$ws.Range('S3').Value = '''11'
$ws.Range('T3').Value = '''1'

# Row 4 - Par
$ws.Range('D4').Value = '''3'
$ws.Range('E4').Value = '''4'
$ws.Range('F4').Value = '''5'
$ws.Range('J4').Value = '''4'
$ws.Range('M4').Value = '''3'
$ws.Range('O4').Value = '''5'
$ws.Range('P4').Value = '''4'
$ws.Range('Q4').Value = '''4'

# Row 5 - Jerome (score)
$ws.Range('B5').Value = '''7'
$ws.Range('C5').Value = '''5'
$ws.Range('D5').Value = '''4'
$ws.Range('E5').Value = '''6'
$ws.Range('G5').Value = '''3'
$ws.Range('H5').Value = '''7'
$ws.Range('I5').Value = '''6'
$ws.Range('J5').Value = '''6'
$ws.Range('K5').Value = '''49'
$ws.Range('L5').Value = '''8'
$ws.Range('M5').Value = '''3'
$ws.Range('N5').Value = '''6'
$ws.Range('O5').Value = '''9'
$ws.Range('P5').Value = '''7'
$ws.Range('R5').Value = '''6'
$ws.Range('S5').Value = '''6'
$ws.Range('U5').Value = '''57'
$ws.Range('V5').Value = '''106'

# Row 6 - Round Score
$ws.Range('B6').Value = '''2'
$ws.Range('C6').Value = '''4'
$ws.Range('D6').Value = '''6'
$ws.Range('F6').Value = '''12'
$ws.Range('G6').Value = '''16'
$ws.Range('I6').Value = '''17'
$ws.Range('J6').Value = '''19'
$ws.Range('K6').Value = '''19'
$ws.Range('M6').Value = '''23'
$ws.Range('N6').Value = '''24'
$ws.Range('O6').Value = '''24'
$ws.Range('P6').Value = '''24'
$ws.Range('Q6').Value = '''26'
$ws.Range('R6').Value = '''27'
$ws.Range('S6').Value = '''28'
$ws.Range('T6').Value = '''30'
$ws.Range('U6').Value = '''11'
$ws.Range('V6').Value = '''30'

# Row 7 - Net Score
$ws.Range('B7').Value = '''2'
$ws.Range('C7').Value = '''2'
$ws.Range('D7').Value = '''2'
$ws.Range('E7').Value = '''2'
$ws.Range('F7').Value = '''4'
$ws.Range('G7').Value = '''4'
$ws.Range('H7').Value = '''0'
$ws.Range('I7').Value = '''1'
$ws.Range('K7').Value = '''19'
$ws.Range('L7').Value = '''1'
$ws.Range('M7').Value = '''3'
$ws.Range('N7').Value = '''1'
$ws.Range('O7').Value = '''0'
$ws.Range('P7').Value = '''0'
$ws.Range('Q7').Value = '''2'
$ws.Range('R7').Value = '''1'
$ws.Range('S7').Value = '''1'
$ws.Range('U7').Value = '''11'
$ws.Range('V7').Value = '''30'

# Row 8 - Tee Club
$ws.Range('D8').Value = '5W'
$ws.Range('E8').Value = '3W'
$ws.Range('H8').Value = '3W'
$ws.Range('I8').Value = '7i'
$ws.Range('M8').Value = 'Pw'
$ws.Range('O8').Value = '1W'
$ws.Range('P8').Value = '5W'
$ws.Range('Q8').Value = '5W'
$ws.Range('R8').Value = '6i'
$ws.Range('S8').Value = '6i'

# Row 9 - Fairways
$ws.Range('K9').Value = '''86%'
$ws.Range('U9').Value = '''43%'
$ws.Range('V9').Value = '''64%'

# Row 10 - GIR
$ws.Range('U10').Value = '''13%'
$ws.Range('V10').Value = '''13%'

# Row 11 - Putts
$ws.Range('B11').Value = '''2'
$ws.Range('C11').Value = '''2'
$ws.Range('D11').Value = '''1'
$ws.Range('E11').Value = '''2'
$ws.Range('F11').Value = '''1'
$ws.Range('G11').Value = '''1'
$ws.Range('H11').Value = '—'
$ws.Range('I11').Value = '''2'
$ws.Range('K11').Value = '''14'
$ws.Range('L11').Value = '''3'
$ws.Range('M11').Value = '''2'
$ws.Range('O11').Value = '—'
$ws.Range('Q11').Value = '''3'
$ws.Range('R11').Value = '''3'
$ws.Range('S11').Value = '''3'
$ws.Range('U11').Value = '''21'
$ws.Range('V11').Value = '''35'

# Row 12 - Sand Shots
$ws.Range('B12').Value = '—'
$ws.Range('C12').Value = '''1'
$ws.Range('D12').Value = '—'
$ws.Range('E12').Value = '''1'
$ws.Range('K12').Value = '''2'
$ws.Range('V12').Value = '''3'

# Row 13 - Penalties
$ws.Range('B13').Value = '—'
$ws.Range('H13').Value = '''1'
$ws.Range('I13').Value = '''1'
$ws.Range('K13').Value = '''2'
$ws.Range('M13').Value = '—'
$ws.Range('O13').Value = '''1'
$ws.Range('T13').Value = '—'
$ws.Range('U13').Value = '''1'
$ws.Range('V13').Value = '''3'
